# The architecture diagram on slide 1 was nudged straight up: every shape's
# vertical position (Top) decreases by 320635 EMU (its horizontal position
# and size are untouched). PowerPoint's COM object model stores Left/Top in
# points (1 pt = 12700 EMU) as a 32-bit Single, and the EMU value written
# back to the OOXML truncates that point value rather than rounding it, so
# we bias each target by half an EMU before converting to points to land on
# the exact integer EMU the diff expects.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$emuPerPoint = 12700
$deltaEmu = 320635

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)

    $currentEmu = [Math]::Round($sh.Top * $emuPerPoint)
    $targetEmu = $currentEmu - $deltaEmu

    $sh.Top = ($targetEmu + 0.5) / $emuPerPoint
}
